$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 414.58823
$ws.Range("I33").Value = 363.2
$ws.Range("K33").Value = 363.2
$ws.Range("M33").Value = -134.2

$ws.Range("H62").Value = 2091.6667
$ws.Range("I62").Value = 2054.3635
$ws.Range("J62").Value = 2150.2856
$ws.Range("K62").Value = 2054.3635
$ws.Range("L62").Value = 2150.2856
$ws.Range("M62").Value = -1430.3635
$ws.Range("N62").Value = -3398.2856

$ws.Range("H65").Value = 2091.6667
$ws.Range("I65").Value = 2054.3635
$ws.Range("J65").Value = 2150.2856
$ws.Range("K65").Value = 10271.8175
$ws.Range("L65").Value = 10751.428
$ws.Range("M65").Value = -7151.817499999999
$ws.Range("N65").Value = -16991.428

$ws.Range("H113").Value = 22226320
$ws.Range("I113").Value = 76924810
$ws.Range("J113").Value = 5059
$ws.Range("K113").Value = 76924810
$ws.Range("L113").Value = 5059
$ws.Range("M113").Value = -76921556
$ws.Range("N113").Value = -11567

$ws.Range("H129").Value = 848.325
$ws.Range("I129").Value = 499.4
$ws.Range("J129").Value = 898.17145
$ws.Range("K129").Value = 1498.2
$ws.Range("L129").Value = 2694.51435
$ws.Range("M129").Value = 3501.8
$ws.Range("N129").Value = -12694.51435

$ws.Range("H132").Value = 33002.547
$ws.Range("I132").Value = 35936.133
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 107808.399
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -105278.399
$ws.Range("N132").Value = -16060.0001

$ws.Range("H137").Value = 36175.535
$ws.Range("I137").Value = 3310.125
$ws.Range("J137").Value = 73736
$ws.Range("K137").Value = 9930.375
$ws.Range("L137").Value = 221208
$ws.Range("M137").Value = -7380.375
$ws.Range("N137").Value = -226308

$ws.Range("H138").Value = 2711.0293
$ws.Range("J138").Value = 3687.7727
$ws.Range("L138").Value = 11063.3181
$ws.Range("N138").Value = -21343.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34707.574
$ws.Range("I32").Value = 47192.652
$ws.Range("J32").Value = 5991.9
$ws.Range("K32").Value = 47192.652
$ws.Range("L32").Value = 5991.9
$ws.Range("M32").Value = -46905.652
$ws.Range("N32").Value = -6565.9

$ws.Range("H45").Value = 2858.484
$ws.Range("I45").Value = 2129.9
$ws.Range("J45").Value = 3205.4285
$ws.Range("K45").Value = 2129.9
$ws.Range("L45").Value = 3205.4285
$ws.Range("M45").Value = -1752.9
$ws.Range("N45").Value = -3959.4285

$ws.Range("H94").Value = 39000
$ws.Range("J94").Value = 39000
$ws.Range("L94").Value = 39000
$ws.Range("N94").Value = -40802

$ws.Range("H102").Value = 7252.75
$ws.Range("I102").Value = 8000
$ws.Range("K102").Value = 8000
$ws.Range("M102").Value = -6378

$ws.Range("H110").Value = 3127.6155
$ws.Range("I110").Value = 2297.6667
$ws.Range("K110").Value = 2297.6667
$ws.Range("M110").Value = -252.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9993.92
$ws.Range("I31").Value = 15523.034
$ws.Range("J31").Value = 2358.476
$ws.Range("K31").Value = 15523.034
$ws.Range("L31").Value = 2358.476
$ws.Range("M31").Value = -15228.034
$ws.Range("N31").Value = -2948.476

$ws.Range("H34").Value = 9993.92
$ws.Range("I34").Value = 15523.034
$ws.Range("J34").Value = 2358.476
$ws.Range("K34").Value = 15523.034
$ws.Range("L34").Value = 2358.476
$ws.Range("M34").Value = -15321.034
$ws.Range("N34").Value = -2762.476

$ws.Range("H58").Value = 28940.111
$ws.Range("I58").Value = 1230.7059
$ws.Range("J58").Value = 500000
$ws.Range("K58").Value = 1230.7059
$ws.Range("L58").Value = 500000
$ws.Range("M58").Value = -1027.7059
$ws.Range("N58").Value = -500406

$ws.Range("H95").Value = 13500
$ws.Range("J95").Value = 13500
$ws.Range("L95").Value = 13500
$ws.Range("N95").Value = -18992

$ws.Range("H99").Value = 4754.95
$ws.Range("I99").Value = 3709.9
$ws.Range("J99").Value = 5800
$ws.Range("K99").Value = 3709.9
$ws.Range("L99").Value = 5800
$ws.Range("M99").Value = -2211.9
$ws.Range("N99").Value = -8796

$ws.Range("H107").Value = 987.9375
$ws.Range("I107").Value = 1266.091
$ws.Range("K107").Value = 1266.091
$ws.Range("M107").Value = 653.9090000000001

$ws.Range("H122").Value = 2621.1667
$ws.Range("I122").Value = 3378.25
$ws.Range("J122").Value = 1107
$ws.Range("K122").Value = 10134.75
$ws.Range("L122").Value = 3321
$ws.Range("M122").Value = -7684.75
$ws.Range("N122").Value = -8221

$ws.Range("H126").Value = 4754.95
$ws.Range("I126").Value = 3709.9
$ws.Range("J126").Value = 5800
$ws.Range("K126").Value = 11129.7
$ws.Range("L126").Value = 17400
$ws.Range("M126").Value = -8659.700000000001
$ws.Range("N126").Value = -22340

$ws.Range("H132").Value = 21694.902
$ws.Range("I132").Value = 22259.621
$ws.Range("J132").Value = 13506.5
$ws.Range("K132").Value = 66778.863
$ws.Range("L132").Value = 40519.5
$ws.Range("M132").Value = -64248.863
$ws.Range("N132").Value = -45579.5

$ws.Range("H134").Value = 7918.2144
$ws.Range("I134").Value = 835
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 2505
$ws.Range("L134").Value = 300000
$ws.Range("M134").Value = 30
$ws.Range("N134").Value = -305070

$ws.Range("H136").Value = 28940.111
$ws.Range("I136").Value = 1230.7059
$ws.Range("J136").Value = 500000
$ws.Range("K136").Value = 3692.1177
$ws.Range("L136").Value = 1500000
$ws.Range("M136").Value = -1142.1177
$ws.Range("N136").Value = -1505100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 787.1622
$ws.Range("J5").Value = 846.2083
$ws.Range("L5").Value = 2538.6249
$ws.Range("N5").Value = -2762.6249

$ws.Range("H41").Value = 522
$ws.Range("I41").Value = 444
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 1332
$ws.Range("L41").Value = 1800
$ws.Range("M41").Value = -994
$ws.Range("N41").Value = -2476

$ws.Range("H94").Value = 5575.857
$ws.Range("I94").Value = 950
$ws.Range("J94").Value = 6346.8335
$ws.Range("K94").Value = 2850
$ws.Range("L94").Value = 19040.5005
$ws.Range("M94").Value = -2174
$ws.Range("N94").Value = -20392.5005

$ws.Range("H131").Value = 118493.37
$ws.Range("J131").Value = 129053.54
$ws.Range("L131").Value = 387160.62
$ws.Range("N131").Value = -397240.62

$ws.Range("H135").Value = 787.1622
$ws.Range("J135").Value = 846.2083
$ws.Range("L135").Value = 7615.8747
$ws.Range("N135").Value = -12685.8747

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 14256
$ws.Range("I102").Value = 26012
$ws.Range("K102").Value = 26012
$ws.Range("M102").Value = -24390

$ws.Range("H107").Value = 888.4706
$ws.Range("I107").Value = 377.85715
$ws.Range("J107").Value = 1245.9
$ws.Range("K107").Value = 377.85715
$ws.Range("L107").Value = 1245.9
$ws.Range("M107").Value = 1542.14285
$ws.Range("N107").Value = -5085.9

$ws.Range("H113").Value = 4799.857
$ws.Range("I113").Value = 2866.6667
$ws.Range("J113").Value = 6249.75
$ws.Range("K113").Value = 2866.6667
$ws.Range("L113").Value = 6249.75
$ws.Range("M113").Value = -696.6667000000002
$ws.Range("N113").Value = -10589.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2160.6667
$ws.Range("I22").Value = 2800.25
$ws.Range("J22").Value = 881.5
$ws.Range("K22").Value = 2800.25
$ws.Range("L22").Value = 881.5
$ws.Range("M22").Value = -2505.25
$ws.Range("N22").Value = -1471.5

$ws.Range("H27").Value = 2160.6667
$ws.Range("I27").Value = 2800.25
$ws.Range("J27").Value = 881.5
$ws.Range("K27").Value = 2800.25
$ws.Range("L27").Value = 881.5
$ws.Range("M27").Value = -2693.25
$ws.Range("N27").Value = -1095.5

$ws.Range("H40").Value = 105619.18
$ws.Range("J40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("N40").Value = -8272

$ws.Range("H122").Value = 3750.75
$ws.Range("J122").Value = 4186.4287
$ws.Range("L122").Value = 12559.2861
$ws.Range("N122").Value = -17459.2861

$ws.Range("H136").Value = 2874.6667
$ws.Range("I136").Value = 2054.111
$ws.Range("J136").Value = 5336.3335
$ws.Range("K136").Value = 6162.333
$ws.Range("L136").Value = 16009.0005
$ws.Range("M136").Value = -3612.333
$ws.Range("N136").Value = -21109.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1603
$ws.Range("I122").Value = 1446.25
$ws.Range("K122").Value = 4338.75
$ws.Range("M122").Value = -1888.75

$ws.Range("H132").Value = 2003.56
$ws.Range("I132").Value = 1685.9546
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5057.8638
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2527.8638
$ws.Range("N132").Value = -18057.9995
